$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "process"
$ws.Range("D2").Value = "OK"
$ws.Range("D3").Value = "OK"
$ws.Range("D4").Value = "OK"
$ws.Range("D5").Value = "OK"

$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "meas_bmi_id"
$ws.Range("A2").Value = 9996
$ws.Range("A3").Value = 9997
$ws.Range("A4").Value = 9998
$ws.Range("A5").Value = 9999

$ws.Columns("A:A").ColumnWidth = 12

$ws.Range("E5").Font.Name = "Calibri"

$ws.Range("A2:A5").Select()
